$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Center crop the image tensor" -> "Center crop the image tensor." (add trailing period)
$ws.Range("B7").Value = "Center crop the image tensor."

# Update the selected/active cell from B13 to B8
$ws.Range("B8").Select()
